$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "TeamsSimpleVacation" sheet after "TeamsMultiple", seeded
#    from the same data as TeamsSimple but with a 3rd vacation interval
#    (written first so the new shared-string order matches the original
#    author's edit order)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$vac = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$vac.Name = "TeamsSimpleVacation"

$vac.Range("A1").Value = "Team ID"
$vac.Range("B1").Value = "Team Name"
$vac.Range("C1").Value = "Developer Id"
$vac.Range("D1").Value = "Developer Name"
$vac.Range("E1").Value = "Developer Vacation Date Intervals"
$vac.Range("F1").Value = "Daily Work Hours"

$vac.Range("A2").Value = 1
$vac.Range("B2").Value = "Team Alpha"
$vac.Range("C2").Value = 1
$vac.Range("D2").Value = "Alice"
$vac.Range("E2").Value = "2025-03-10;2025-03-15|2025-05-10;2025-05-15|2025-05-19;2025-05-19"
$vac.Range("F2").Value = 6

$vac.Range("A3").Value = 1
$vac.Range("B3").Value = "Team Alpha"
$vac.Range("C3").Value = 2
$vac.Range("D3").Value = "Bob"
$vac.Range("F3").Value = 7

$vac.Range("A4").Value = 1
$vac.Range("B4").Value = "Team Alpha"
$vac.Range("C4").Value = 3
$vac.Range("D4").Value = "Charlie"
$vac.Range("E4").Value = "2025-04-05;2025-04-10|"
$vac.Range("F4").Value = 7

$vac.Range("A5").Value = 1
$vac.Range("B5").Value = "Team Alpha"
$vac.Range("C5").Value = 4
$vac.Range("D5").Value = "Dave"
$vac.Range("E5").Value = "2025-02-25;2025-02-28|"
$vac.Range("F5").Value = 6

# column widths to roughly match the source sheet
$vac.Columns.Item(1).ColumnWidth = 12.33
$vac.Columns.Item(2).ColumnWidth = 12.33
$vac.Columns.Item(3).ColumnWidth = 16.55
$vac.Columns.Item(4).ColumnWidth = 16.44
$vac.Columns.Item(5).ColumnWidth = 15.5

# turn the data range into a table, like the other sheets
$tbl = $vac.ListObjects.Add(1, $vac.Range("A1:F5"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table14"

$vac.Activate() | Out-Null
$vac.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. TeamsSimple: tweak one developer's vacation intervals + an hours value
# ---------------------------------------------------------------------------
$simple = $wb.Worksheets.Item("TeamsSimple")
$simple.Range("E2").Value = "2025-03-10;2025-03-15|2025-05-10;2025-05-15|"
$simple.Range("F3").Value = 7

# ---------------------------------------------------------------------------
# 3. Restore the selection on TeamsSimple to what the diff records
# ---------------------------------------------------------------------------
$simple.Activate() | Out-Null
$simple.Range("E2").Select() | Out-Null
